# Adds 7 new firearm-type weapons into the WeaponStats table, keeping the
# existing (alphabetically sorted) columns intact by inserting a fresh
# column at each new weapon's correct alphabetical slot and filling in its
# five rows of data (Name / Damage / Ranged? / Martial? / Properties).
#
# Insertion is done left-to-right using the *final* target column letters;
# since each insert shifts everything at/after that column one slot to the
# right, performing the inserts in ascending column order reproduces the
# final table layout exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newWeapons = @(
    @{ Col = "A";  Name = "Bad News";    Damage = "2d12=P"; Ranged = $true;  Martial = $false; Props = "Ammunition=Two-handed=Reload 1=Misfire 3=Range (200/800)" },
    @{ Col = "C";  Name = "Blunderbuss"; Damage = "2d8=P";  Ranged = $true;  Martial = $false; Props = "Ammunition=Reload 1=Misfire 2=Range (15/60)" },
    @{ Col = "R";  Name = "Hand-Mortar"; Damage = "2d8=F";  Ranged = $true;  Martial = $false; Props = "Ammunition=Reload 1=Misfire 3=Explosive=Ranged (30/600)" },
    @{ Col = "AA"; Name = "Musket";      Damage = "1d12=P"; Ranged = $true;  Martial = $false; Props = "Ammunition=Ranged (120/480)=Two-handed=Reload 1=Misfire 2" },
    @{ Col = "AC"; Name = "Palm Pistol"; Damage = "1d8=P";  Ranged = $true;  Martial = $false; Props = "Ammunition=Ranged (40/160)=Light=Reload 1=Misfire 1" },
    @{ Col = "AD"; Name = "Pepperbox";   Damage = "1d10=P"; Ranged = $true;  Martial = $false; Props = "Ammunition=Ranged (80/320)=Reload 6=Misfire 2" },
    @{ Col = "AF"; Name = "Pistol";      Damage = "1d10=P"; Ranged = $true;  Martial = $false; Props = "Ammunition=Ranged (60/240)=Reload 4=Misfire 1" }
)

foreach ($weapon in $newWeapons) {
    $col = $weapon.Col

    $ws.Columns($col).Insert()

    if ($col -eq "A") {
        # Column A has no left neighbour to inherit formatting from on
        # insert, so explicitly copy the style from its (now pushed-over)
        # neighbour column B to match the rest of the table (style index 1).
        $ws.Range("B1:B5").Copy()
        $ws.Range("A1:A5").PasteSpecial(-4122)
    }

    $ws.Range($col + "1").Value = $weapon.Name
    $ws.Range($col + "2").Value = $weapon.Damage
    $ws.Range($col + "3").Value = if ($weapon.Ranged) { "true=" } else { "false=" }
    $ws.Range($col + "4").Value = if ($weapon.Martial) { "true=" } else { "false=" }
    $ws.Range($col + "5").Value = $weapon.Props
}
